$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Retrospectiva")

$ws.Range("A2").Value = "Github anduvo bien toda la semana hasta el ultimo momento"
$ws.Range("C3").Value = "primero hacer commit y dspues sincronizar"
$ws.Range("C4").Value = "la idea es no tocar lo mismo"
$ws.Range("C2").Clear()
$ws.Range("C2").Value = "mejorar github"

$ws.Range("C2").Select()
